$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the stored precision of A68 (same underlying value, extra decimal digit)
$ws.Range("A68").Value = 44381.76711011458

# Append the new data row retrieved on 2021-07-05 (row 69)
$ws.Range("A69").Value = 44382.76768740083
$ws.Range("B69").Value = 78624
$ws.Range("C69").Value = 66315
$ws.Range("D69").Value = 3569
$ws.Range("E69").Value = 2138
$ws.Range("F69").Value = 1524
$ws.Range("G69").Value = 20868
$ws.Range("H69").Value = 1536
$ws.Range("I69").Value = 878
$ws.Range("J69").Value = 195
